# Toronto roster: swap the full row content (No., Player, Pos, Ht, Wt,
# Birth Date, Exp, College, bbref url) between rows 15/16 and between
# rows 17/18.  Row numbers (column A) stay put; only what is displayed
# in each row moves, exactly like the source diff shows (shared-string
# text is relabeled while the existing hyperlink relationships for
# column K are left completely untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell as TEXT (not a
# number), without introducing a new cell style (quote-prefix would add
# one). We stage the text via a formula that evaluates to a string in a
# scratch cell, copy it, and paste-special just the value into the
# destination, which preserves the shared-string (text) type cleanly.
function Set-TextValue($rangeAddr, $text) {
    $ws.Range("M1").Formula = "=""$text"""
    $ws.Range("M1").Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
    $ws.Range("M1").Clear()
}

# --- Row 15 becomes what Row 16 used to show -------------------------
$ws.Range("B15").Value2 = 1
$ws.Range("C15").Value2 = "Will Barton"
$ws.Range("D15").Value2 = "SG"
$ws.Range("E15").Value2 = "6-6"
$ws.Range("F15").Value2 = 181
$ws.Range("G15").Value2 = "January 6, 1991"
Set-TextValue "I15" "10"
$ws.Range("J15").Value2 = "Memphis"
$ws.Range("K15").Value2 = "https://www.basketball-reference.com/players/b/bartowi01.html"

# --- Row 16 becomes what Row 15 used to show -------------------------
$ws.Range("B16").Value2 = 32
$ws.Range("C16").Value2 = "Otto Porter Jr."
$ws.Range("D16").Value2 = "SF"
$ws.Range("E16").Value2 = "6-8"
$ws.Range("F16").Value2 = 198
$ws.Range("G16").Value2 = "June 3, 1993"
Set-TextValue "I16" "9"
$ws.Range("J16").Value2 = "Georgetown"
$ws.Range("K16").Value2 = "https://www.basketball-reference.com/players/p/porteot01.html"

# --- Row 17 becomes what Row 18 used to show -------------------------
$ws.Range("B17").Value2 = 11
$ws.Range("C17").Value2 = "Joe Wieskamp"
$ws.Range("D17").Value2 = "SF"
$ws.Range("E17").Value2 = "6-6"
$ws.Range("F17").Value2 = 212
$ws.Range("G17").Value2 = "August 23, 1999"
Set-TextValue "I17" "1"
$ws.Range("J17").Value2 = "Iowa"
$ws.Range("K17").Value2 = "https://www.basketball-reference.com/players/w/wieskjo01.html"

# --- Row 18 becomes what Row 17 used to show -------------------------
$ws.Range("B18").Value2 = 8
$ws.Range("C18").Value2 = "Ron Harper Jr. (TW)"
$ws.Range("D18").Value2 = "SF"
$ws.Range("E18").Value2 = "6-6"
$ws.Range("F18").Value2 = 245
$ws.Range("G18").Value2 = "April 12, 2000"
$ws.Range("I18").Value2 = "R"
$ws.Range("J18").Value2 = "Rutgers University"
$ws.Range("K18").Value2 = "https://www.basketball-reference.com/players/h/harpero02.html"
